# Updated latest Guinea master data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# New column order: word, descr, lang_code, is_active, cr_by, cr_dtimes,
# upd_by, upd_dtimes, is_deleted, del_dtimes
$ws.Cells.Item(1, 1).Value = "word"
$ws.Cells.Item(1, 2).Value = "descr"
$ws.Cells.Item(1, 3).Value = "lang_code"
$ws.Cells.Item(1, 4).Value = "is_active"
$ws.Cells.Item(1, 5).Value = "cr_by"
$ws.Cells.Item(1, 6).Value = "cr_dtimes"
$ws.Cells.Item(1, 7).Value = "upd_by"
$ws.Cells.Item(1, 8).Value = "upd_dtimes"
$ws.Cells.Item(1, 9).Value = "is_deleted"
$ws.Cells.Item(1, 10).Value = "del_dtimes"

# --- Data rows ----------------------------------------------------------
$descr   = "Mot sur la liste noire"
$lang    = "fra"
$crBy    = "superadmin"
$crDate  = 45079.576874606479
$nullStr = "NULL"

$words = @("Merde", "pute", "putain", "bon sang", "chier", "imbÃ©cile", "faire foutre")

for ($i = 0; $i -lt $words.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $words[$i]
    $ws.Cells.Item($r, 2).Value = $descr
    $ws.Cells.Item($r, 3).Value = $lang
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = $crBy
    $ws.Cells.Item($r, 6).Value = $crDate
    $ws.Cells.Item($r, 7).Value = $nullStr
    $ws.Cells.Item($r, 8).Value = $nullStr
    $ws.Cells.Item($r, 9).Value = $false
    $ws.Cells.Item($r, 10).Value = $nullStr
}

# cr_dtimes column uses the built-in "mm:ss.0" number format (numFmtId 47)
$ws.Range("F2:F8").NumberFormat = "mm:ss.0"

# --- Selection (matches saved cursor position in the target file) -------
[void]$ws.Range("D13").Select()
